# Update the monthly expense values on every month sheet (January - December).
# Per commit message: replace the placeholder test data with real current/max
# values - Gas: Current 90 / Max 150, Rent: Current 0 / Max 0, Food: Current 0 / Max 0.
# Column B holds the "Current" value, column C holds the "Max" value.
# A leading apostrophe is used so the numeric-looking text stays stored as text
# (matching the original inline-string cell contents) instead of being
# auto-converted to a number by Excel.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Gas (row 2): Current = 90, Max = 150
    $ws.Range("B2").Value = "'90"
    $ws.Range("C2").Value = "'150"

    # Rent (row 3): Current = 0, Max = 0
    $ws.Range("B3").Value = "'0"
    $ws.Range("C3").Value = "'0"

    # Food (row 4): Current = 0, Max = 0
    $ws.Range("B4").Value = "'0"
    $ws.Range("C4").Value = "'0"
}
